# ApachePOI Citizen functionality cozumu
# Updates the "testCitizen" sheet's string data (ulais114x / urbsXX values),
# tweaks the sheet view (zoom + selection) and the column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")
[void]$ws.Activate()

# --- Shared-string content edits (column A / column B values) ---
$ws.Range("A1").Value = "ulais11451"
$ws.Range("A2").Value = "ulais11462"
$ws.Range("A3").Value = "ulais11473"
$ws.Range("A4").Value = "ulais11483"
$ws.Range("A5").Value = "ulais11494"
$ws.Range("A6").Value = "ulais11505"
$ws.Range("A7").Value = "ulais11516"
$ws.Range("A8").Value = "ulais11527"

$ws.Range("B1").Value = "urbs131"
$ws.Range("B2").Value = "urbs141"
$ws.Range("B3").Value = "urbs151"
$ws.Range("B4").Value = "urbs161"
$ws.Range("B5").Value = "urbs171"
$ws.Range("B6").Value = "urbs181"
$ws.Range("B7").Value = "urbs191"
$ws.Range("B8").Value = "urbs201"

# --- Column A width: widen very slightly (20.285... -> 20.332...) ---
$ws.Columns.Item(1).ColumnWidth = 19.5

# --- View changes: zoom in to 160% and move the selection to B10 ---
$excel.ActiveWindow.Zoom = 160
[void]$ws.Range("B10").Select()
